# Adds "subrogation expense" (cause of loss / subrogate loss) rows to the
# Navigation sheet: 5 new rows (84-88) in columns D/E with new key/xpath
# string pairs, mirroring the existing Subrogate* rows pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("SubrogateLossRadioBox", "/html/body/div[2]/div/form/div/div[2]/div[1]/div/div[1]/div/div[2]/div/div[1]/div/div/div[3]/div/div/div/div[6]/div/div/div[2]/div[1]/div/table/tbody/tr[4]/td[1]/div/input"),
    @("SubrogateLossSubmit", "/html/body/div[2]/div/form/div/div[2]/div[1]/div/div[1]/div/div[2]/div/div[1]/div/div/div[3]/div/div/div/div[6]/div/div/div[2]/div[2]/button[1]"),
    @("SubrogateLossUpdate", "/html/body/div[2]/div/form/div/div[2]/div[1]/div/div[1]/div/div[2]/div/div[1]/div/div/div[3]/div/div/div/div[1]/div/div/div[2]/div[2]/div/table/tbody/tr[2]/td[9]/a"),
    @("SubrogateLossAmtUpdate", "/html/body/div[2]/div/form/div/div[2]/div[1]/div/div[1]/div/div[2]/div/div[1]/div/div/div[3]/div/div/div/div[4]/div/div/div[2]/div[1]/div[2]/div[1]/div[4]/div/div/input"),
    @("SubrogateLossAmtSubmit", "/html/body/div[2]/div/form/div/div[2]/div[1]/div/div[1]/div/div[2]/div/div[1]/div/div/div[3]/div/div/div/div[4]/div/div/div[2]/div[2]/button[1]")
)

$startRow = 84
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $pair = $newRows[$i]
    $ws.Range("D$r").Value = $pair[0]
    $ws.Range("E$r").Value = $pair[1]
}

# Move the selection/active cell to mirror the author's final cursor
# position after adding the new rows.
$ws.Range("D92").Select()
